{"js": "// Color the text of the \"Create new Unit Test Project and add reference to\n// the \"CustomLinkedList\".\" list item red (FF0000), matching the author's\n// edit. We locate the paragraph by its distinctive text and set the font\n// color on the whole paragraph range (covers every run + the paragraph\n// mark's run properties).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"Create new Unit Test Project and add reference to the\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.indexOf(needle) !== -1) {\n    para.font.color = \"#FF0000\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Color the text of the \"Create new Unit Test Project and add reference to\n# the \"CustomLinkedList\".\" list item red (FF0000), matching the author's\n# edit. We locate the paragraph by its distinctive text and set the font\n# color on the whole paragraph range (covers every run + the paragraph\n# mark's run properties).\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -like \"Create new Unit Test Project and add reference to the*CustomLinkedList*\") {\n        $p.Range.Font.Color = [Microsoft.Office.Interop.Word.WdColor]::wdColorRed\n    }\n}\n"}
